$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C3").Value = 0.606439356064394
$ws.Range("C7").Value = 0.22027797220278
$ws.Range("C8").Value = 0.78002199780022

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("C3").Value = 0.917908209179082
$ws.Range("C6").Value = 0.26037396260374
$ws.Range("C7").Value = 0.863813618638136

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C3").Value = 0.2004799520048
$ws.Range("C6").Value = 0.538446155384462
$ws.Range("C7").Value = 0.789221077892211
$ws.Range("C8").Value = 0.695630436956304

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C3").Value = 0.780621937806219
$ws.Range("C6").Value = 0.475852414758524
$ws.Range("C7").Value = 0.304069593040696
$ws.Range("C8").Value = 0.833216678332167

$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("C3").Value = 0.806119388061194
$ws.Range("C5").Value = 0.502549745025497
$ws.Range("C6").Value = 0.934406559344066
$ws.Range("C7").Value = 0.847915208479152

$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("C3").Value = 0.126187381261874
$ws.Range("C6").Value = 0.0217978202179782
$ws.Range("C7").Value = 0.848115188481152
$ws.Range("C8").Value = 0.485651434856514

$ws = $wb.Worksheets.Item("Germ cell tumor")
$ws.Range("C3").Value = 0.623437656234377
$ws.Range("C5").Value = 0.0085991400859914
$ws.Range("C6").Value = 0.852514748525148
$ws.Range("C7").Value = 0.365863413658634
$ws.Range("C8").Value = 0.269273072692731

$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("C3").Value = 0.740925907409259
$ws.Range("C5").Value = 0.004999500049995
$ws.Range("C6").Value = 0.16988301169883
$ws.Range("C7").Value = 0.570642935706429

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("C3").Value = 0.0491950804919508
$ws.Range("C5").Value = 0.0003999600039996
$ws.Range("C6").Value = 0.87971202879712
$ws.Range("C7").Value = 0.992800719928007

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("C3").Value = 0.395560443955604
$ws.Range("C6").Value = 0.263373662633737
$ws.Range("C7").Value = 0.279372062793721
$ws.Range("C8").Value = 0.966803319668033

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("C2").Value = 0.652232662013525
$ws.Range("C3").Value = 0.618438156184382
$ws.Range("C6").Value = 0.017998200179982
$ws.Range("C7").Value = 0.777622237776222

$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("C3").Value = 0.367763223677632
$ws.Range("C5").Value = 0.0003999600039996
$ws.Range("C6").Value = 0.514048595140486
$ws.Range("C7").Value = 0.273872612738726
$ws.Range("C8").Value = 0.17048295170483

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("C2").Value = 0.957322984148866
$ws.Range("C3").Value = 0.568943105689431
$ws.Range("C6").Value = 0.803219678032197
$ws.Range("C7").Value = 0.012998700129987
$ws.Range("C8").Value = 0.661933806619338

$ws = $wb.Worksheets.Item("Meningioma")
$ws.Range("C3").Value = 0.93950604939506
$ws.Range("C5").Value = 0.0007999200079992
$ws.Range("C6").Value = 0.683131686831317
$ws.Range("C7").Value = 0.318868113188681

$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("C4").Value = 0.0003999600039996
$ws.Range("C5").Value = 0.145385461453855
$ws.Range("C6").Value = 0.0598940105989401
$ws.Range("C7").Value = 0.304569543045695

$ws = $wb.Worksheets.Item("Oligodendroglioma")
$ws.Range("C3").Value = 0.168683131686831
$ws.Range("C5").Value = 0.21037896210379
$ws.Range("C7").Value = 0.215478452154785
